$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new "Write Latency / average" (column Q) value.
# Values are entered with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr/text cells) instead of auto-converting
# the numeric-looking strings into real numbers.
$newQValues = @{
    3  = "45068.10"
    4  = "50624.71"
    5  = "349783.45"
    6  = "7619.81"
    7  = "156258.79"
    8  = "30571.82"
    9  = "9532.50"
    10 = "124094.45"
    11 = "3971.32"
    12 = "20.52"
    13 = "243746.25"
    14 = "19.71"
    15 = "1524.04"
    16 = "1723.49"
    17 = "27820.85"
    18 = "5828.10"
    19 = "15084.99"
    20 = "18.75"
    21 = "1742.68"
    22 = "8301.50"
    23 = "13573.38"
}

foreach ($row in $newQValues.Keys) {
    $ws.Range("Q$row").Value = "'" + $newQValues[$row]
}
